# Daily attendance processing - reorder "Recorded By" (column G) entries
# so that the literal "System" token moves from the front of the
# comma-separated list to the end, for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# Rows whose "Recorded By" value begins with "System, " (i.e. has other
# recorders alongside it) need System moved to the end of the list.
$rows = @(2,3,4,5,6,7,8,28,29,30,31,32,33,34,54,55,56,57,58,59,60,80,81,82,106,107,108,132,133,134)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and $value -like "System,*") {
        $parts = $value -split ",\s*"
        $newParts = @()
        $removed = $false
        foreach ($p in $parts) {
            if (-not $removed -and $p -eq "System") {
                $removed = $true
                continue
            }
            $newParts += $p
        }
        $newParts += "System"
        $cell.Value = [string]::Join(", ", $newParts)
    }
}
